$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing data (and column
# widths) from A:E to B:F.
$ws.Columns("A:A").Insert()

# Populate the new column A with the "last update" data.
$ws.Range("A1").Value = "update"
$ws.Range("A2").Value = 20150809
$ws.Range("A3").Value = 20150809

# The column insert does not move the hyperlinks that lived on the old
# B3:E3 (now C3:F3), so remove them and re-add them on the shifted cells.
$ws.Range("B3").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("C3"), "http://pps.sagepub.com/content/by/year/")
$ws.Hyperlinks.Add($ws.Range("D3"), "http://pps.sagepub.com/content/by/year/[0-9]{4}")
$ws.Hyperlinks.Add($ws.Range("E3"), "http://pps.sagepub.com/content/vol[0-9]{1,}/issue[0-9]{1,}/")
$ws.Hyperlinks.Add($ws.Range("F3"), "http://pps.sagepub.com/content/[0-9]{1,}/[0-9]{1,}/[0-9]{1,}.abstract")

# Re-apply the original Hyperlink cell style so the re-added links use the
# same style as before instead of a freshly duplicated one.
$ws.Range("C3:F3").Style = "Hyperlink"

# Update the active selection to match the new layout.
[void]$ws.Range("A4").Select()
